$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.007649031253165467, 0.7060764589429129, 1.09292080587746, 1.045428527388391, 1.055802796057253, 51)
    3  = @(0.1975818930959411, 0.6676882907705637, 0.9882928130720809, 0.9941291732325739, 0.9841883678757558, 50)
    4  = @(0.01857348150839668, 0.6965285593079298, 1.031947222433671, 1.015848031170839, 1.026203664836019, 49)
    5  = @(0.197272805308272, 0.6719669695526704, 1.016022479395869, 1.007979404251828, 0.9989471494200659, 48)
    6  = @(0.008725544551307194, 0.6890175284598777, 1.005961649362776, 1.002976395217144, 1.01378133274956, 47)
    7  = @(0.1656324924955881, 0.7121334300484996, 1.009049436880241, 1.004514527958775, 1.001713026750313, 46)
    8  = @(-0.03355332013530137, 0.6655533201353013, 0.9216836371613637, 0.9600435600332746, 0.9702987067500585, 45)
    9  = @(0.1120563280973417, 0.6615800355390219, 0.9519011468707899, 0.9756542148070646, 0.980402849722379, 44)
    10 = @(0.02418838590491716, 0.6296238101770754, 0.8605132227654934, 0.9276385194489788, 0.9382977269448005, 43)
    11 = @(0.09524948657816748, 0.6695304106776226, 0.9431468323975638, 0.9711574704431634, 0.9781905077005605, 42)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
